$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-organize rows 49-58 ---
# Before state:
#   Row 49: tooltip.SetActiveAtStartToggle (ht=30)
#   Rows 50-58: pause.* translation rows
# Target state:
#   Rows 49-57: pause.* translation rows (shifted up by one)
#   Row 58: selection.NoObjectSelected (new)
#   Row 59: selection.MultipleObjectsSelected (new)
#   Row 60: tooltip.SetActiveAtStartToggle (moved to bottom, ht=30)

# Remove the tooltip row from its current position (row 49); this shifts
# rows 50-58 up to become rows 49-57.
$ws.Rows.Item(49).Delete()

# Add the two new "selected object" header translation rows at 58 and 59.
$ws.Range("A58").Value = "selection.NoObjectSelected"
$ws.Range("B58").Value = "No Object Selected"
$ws.Range("C58").Value = "Sin Objeto Seleccionado"

$ws.Range("A59").Value = "selection.MultipleObjectsSelected"
$ws.Range("B59").Value = "Multiple Objects Selected"
$ws.Range("C59").Value = "Multiples Objetos Seleccionados"

# Re-append the tooltip row at the very bottom (row 60), preserving its
# original content and its taller (30pt) row height / left-aligned English
# tooltip text.
$ws.Range("A60").Value = "tooltip.SetActiveAtStartToggle"
$ws.Range("B60").Value = "Defines if the object will be [00FFFF]enabled[-] or [FF0000]not[-] at the beginning of the level."
$ws.Range("C60").Value = "Define si el objeto estará [00FFFF]activado[-] o [FF0000]no[-] al inicio del nivel."
$ws.Range("B60").HorizontalAlignment = -4131
$ws.Rows.Item(60).RowHeight = 30

# pause.ExitPopup.Content row (now row 54) keeps its taller (45pt) wrapped row height.
$ws.Rows.Item(54).RowHeight = 45

# --- Update the visible view/selection to match the edited workbook ---
# (the headless runtime does not persist window scroll position / topLeftCell,
# so only the active cell / selection is set here)
$ws.Range("D60").Select()
